$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column for all existing data rows (2-501)
#    from 45205 to 45206.
$ws.Range("C2:C501").Value = 45206

# 2. Ensure row 501 has an explicit custom row height (ht="15" customHeight="1").
$ws.Rows.Item(501).RowHeight = 15

# 3. Add the new record as row 502.
$ws.Range("A502").Value = "A 48367-2023"
$ws.Range("B502").Value = 45205
$ws.Range("C502").Value = 45206
$ws.Range("D502").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E502").Value = "VÄSTERÅS"
$ws.Range("G502").Value = 5.9
$ws.Range("H502").Value = 0
$ws.Range("I502").Value = 0
$ws.Range("J502").Value = 0
$ws.Range("K502").Value = 0
$ws.Range("L502").Value = 0
$ws.Range("M502").Value = 0
$ws.Range("N502").Value = 0
$ws.Range("O502").Value = 0
$ws.Range("P502").Value = 0
$ws.Range("Q502").Value = 0
$ws.Range("R502").Value = ""

# Copy cell formatting (number formats / wrap text) from row 501 onto row 502
# so the date columns (B, C) keep their date format and R keeps wrap text,
# matching the style used throughout the sheet.
$ws.Range("A501:E501").Copy()
$ws.Range("A502:E502").PasteSpecial(-4122)
$ws.Range("G501:R501").Copy()
$ws.Range("G502:R502").PasteSpecial(-4122)
